$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - reorder block labels
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Update data rows 2-7 with the new block order values
$values = @(
    @(0, 0, 0, 0, 1, 0),
    @(0, 0, 1, 0, 0, 0),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(0, 1, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowValues[$j]
    }
}
